# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled
# update). Every price (column D) and 1h-volume-change (column E) cell is
# re-stamped with the latest scraped figures; rows 13/14 additionally swap
# rank between "Avalanche" and "TRON" (coin, link, price and change all move
# together), matching the upstream data re-fetch that produced this diff.
#
# Price cells that look like plain numbers (e.g. "0.998", "576.42") are
# written with a leading apostrophe so Excel keeps them as text, exactly
# like the source data (prices with "." as a thousands separator, e.g.
# "61.605.70", already fail numeric parsing and stay text on their own).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.605.70'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '3.393.56'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''576.42'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').Value = '''140.94'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').Value = '''0.389'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').Value = '3.972.69'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.125'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '''28.44'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '3.388.90'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = '''0.0000171'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '61.547.06'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').Value = '''13.69'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').Value = '''9.00'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').Value = '''390.13'
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('D22').Value = '''75.11'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '''0.0000113'
$ws.Range('E25').Value = '  -4.42%  '
$ws.Range('E26').Value = '  +6.23%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '''7.30'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').Value = '''8.06'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '''1.38'
$ws.Range('E32').Value = '  -3.81%  '
$ws.Range('D33').Value = '''23.52'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('D35').Value = '''167.69'
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('E36').Value = '  +0.56%  '
$ws.Range('D37').Value = '3.427.50'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = '''1.48'
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').Value = '''0.0772'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  -9.34%  '
$ws.Range('D41').Value = '''0.779'
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('D42').Value = '''4.46'
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('D43').Value = '''1.67'
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('D45').Value = '2.457.51'
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('D46').Value = '''6.71'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('D47').Value = '''22.73'
$ws.Range('E47').Value = '  -2.81%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  -3.13%  '
$ws.Range('D50').Value = '''2.05'
$ws.Range('E50').Value = '  -2.25%  '
$ws.Range('D51').Value = '''0.207'
$ws.Range('E51').Value = '  -1.55%  '
